$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ----------------------------------------------------------------------
# New column G: "colour code" header + per-row hex colour values
# ----------------------------------------------------------------------

# Header cell - reuse the same look as the rest of row 1 (F1's style)
$ws.Range("G1").Value = "colour code"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Plain-value rows (no special font / alignment applied)
$plainValues = @{
    9  = "#ab9b93"
    10 = "#a78879"
    12 = "#d2d5f4"
    13 = "#d2e9f4"
    15 = "#dadbec"
    16 = "#ffffff"
    17 = "#fbf0f7"
    18 = "#572a46"
    19 = "#d2bf56"
    20 = "#7e82be"
    21 = "#e9dfab"
    22 = "#c45e9e"
    23 = "#4dbd5e"
    24 = "#c395e5"
}
foreach ($row in $plainValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $plainValues[$row]
}

# Rows styled with a small black "Lucida Console" font, left/centre aligned
$lucidaRows = @{
    2  = "#fddecf"
    6  = "#474f67"
    7  = "#f4ddd2"
    8  = "#696d79"
    11 = "#989cb4"
    14 = "#c6c9c5"
}
foreach ($row in $lucidaRows.Keys) {
    $ws.Cells.Item($row, 7).Value = $lucidaRows[$row]
}

$g2 = $ws.Cells.Item(2, 7)
$g2.Font.Name = "Lucida Console"
$g2.Font.Size = 7
$g2.Font.Color = 0
$g2.Font.Family = 3
$g2.HorizontalAlignment = -4131   # xlLeft
$g2.VerticalAlignment = -4108     # xlCenter

$g2.Copy()
foreach ($row in @(6, 7, 8, 11, 14)) {
    $ws.Cells.Item($row, 7).PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = $false

# Rows styled with a small "Consolas" font (orange-ish colour), centred vertically
$consolasRows = @{
    3 = "#FDDECF"
    4 = "#384f57"
    5 = "#ec8f6e"
}
foreach ($row in $consolasRows.Keys) {
    $ws.Cells.Item($row, 7).Value = $consolasRows[$row]
}

$g3 = $ws.Cells.Item(3, 7)
$g3.Font.Name = "Consolas"
$g3.Font.Size = 7
$g3.Font.Color = 7901646
$g3.Font.Family = 3
$g3.VerticalAlignment = -4108     # xlCenter

$g3.Copy()
foreach ($row in @(4, 5)) {
    $ws.Cells.Item($row, 7).PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = $false

# ----------------------------------------------------------------------
# Sheet view: selection moved to G23 (scrolled so row 8 is visible)
# ----------------------------------------------------------------------
$ws.Range("G23").Select()
